$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet currently has 4 item rows (4..7), a totals row (8) and a
# footer row (9). The new version adds 3 more item rows (new rows
# 8,9,10), pushing the totals row to 11 and the footer row to 12.
# ------------------------------------------------------------------

# Insert three blank rows before the current totals row (row 8), one
# at a time so each new row lands at position 8 and pushes the rest
# down.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# Clone the formatting (styles + merges layout source) of row 7 (a
# normal item row) onto the three freshly inserted rows so they pick
# up the exact same cell styles used by the other item rows.
$ws.Range("A7:N7").Copy()
$ws.Range("A8:N10").PasteSpecial(-4122)

# Match the row heights seen in the target workbook.
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 24.75
$ws.Rows.Item(10).RowHeight = 25.5

# Re-create the merged cell layout (B:G, H:K, L:M) for each new row,
# matching the layout used by the existing item rows.
$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()

$ws.Range("B9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()

$ws.Range("B10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()

# ------------------------------------------------------------------
# Update the existing item rows (4..7) with their new figures.
# ------------------------------------------------------------------
$ws.Range("B4").Value = "TELFAST 180MG 20 F.C. TABS"
$ws.Range("H4").Value = "1:0"
$ws.Range("L4").Value = 52
$ws.Range("N4").Value = "0:2"

$ws.Range("B5").Value = "URSOFALK 250MG 20 CAPS."
$ws.Range("H5").Value = "1:0"
$ws.Range("L5").Value = 80
$ws.Range("N5").Value = "0:2"

$ws.Range("B6").Value = "WELLMETAZONE 0.1% CREAM 40 GM"
$ws.Range("H6").Value = "0:0"
$ws.Range("L6").Value = 122
$ws.Range("N6").Value = "1:0"

$ws.Range("B7").Value = "كريم فاتيكا 125 مل"
$ws.Range("H7").Value = "2:0"
$ws.Range("L7").Value = 26
$ws.Range("N7").Value = "1:0"

# ------------------------------------------------------------------
# Fill in the three new item rows.
# ------------------------------------------------------------------
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "WELLMETAZONE 0.1% CREAM 40 GM"
$ws.Range("H8").Value = "0:0"
$ws.Range("L8").Value = 56
$ws.Range("N8").Value = "1:0"

$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "سرنجات 5 سم"
$ws.Range("H9").Value = "-1:0"
$ws.Range("L9").Value = 2
$ws.Range("N9").Value = "1:0"

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "كريم فاتيكا 125 مل"
$ws.Range("H10").Value = "2:0"
$ws.Range("L10").Value = 50
$ws.Range("N10").Value = "1:0"

# ------------------------------------------------------------------
# Update the grand-total cell (now on row 11).
# ------------------------------------------------------------------
$ws.Range("K11").Value = 388

# The footer row (now row 12) has a slightly different height in the
# target workbook than the original row 9 it was shifted down from.
$ws.Rows.Item(12).RowHeight = 16.5
